$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, border, centered alignment) from H1 header
# cell so the new header cells I1/J1 match the style of the existing ones.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# PasteSpecial for formats can also paste values in some situations, so
# re-assert the text values after the paste just to be safe.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$excel.CutCopyMode = 0

# --- Data rows 2-55: new columns I (I0) and J (IF) ---
$data = @"
2,7,8
3,8,9
4,9,9
5,6,6
6,7,7
7,6,6
8,5,6
9,8,8
10,7,7
11,5,5
12,5,5
13,6,7
14,6,6
15,8,8
16,8,8
17,8,8
18,8,8
19,5,6
20,5,6
21,8,8
22,5,5
23,6,6
24,8,8
25,9,9
26,6,6
27,7,7
28,5,5
29,5,5
30,4,4
31,5,6
32,8,9
33,6,6
34,6,6
35,10,10
36,4,5
37,5,6
38,7,8
39,8,9
40,4,5
41,2,3
42,7,7
43,6,7
44,5,5
45,8,8
46,8,8
47,9,9
48,7,7
49,8,8
50,6,6
51,8,9
52,6,7
53,8,8
54,6,6
55,7,7
"@

$rows = $data -split "`n"
foreach ($line in $rows) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $iVal = [int]$parts[1]
    $jVal = [int]$parts[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Host "Applied I0/IF columns"
